$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24
$ws.Range("D24").Value = 44468
$ws.Range("L24").Value = "Primera"
$ws.Range("M24").Value = 300
$ws.Range("N24").Value = 11000
$ws.Range("O24").Value = 12000
$ws.Range("P24").Value = 11500
$ws.Range("S24").Value = 575

# Row 25
$ws.Range("D25").Value = 44399
$ws.Range("K25").Value = "Clemenuless"
$ws.Range("N25").Value = 13000
$ws.Range("O25").Value = 14000
$ws.Range("P25").Value = 13500
$ws.Range("R25").Value = "Región Metropolitana"
$ws.Range("S25").Value = 675

# Row 26
$ws.Range("D26").Value = 44258
$ws.Range("K26").Value = "Murcott"
$ws.Range("L26").Value = "Segunda"
$ws.Range("M26").Value = 400
$ws.Range("N26").Value = 17000
$ws.Range("O26").Value = 18000
$ws.Range("P26").Value = 17500
$ws.Range("R26").Value = "Región de Coquimbo"
$ws.Range("S26").Value = 875

# Row 27
$ws.Range("D27").Value = 44321
$ws.Range("K27").Value = "Clementina"
$ws.Range("L27").Value = "Primera"
$ws.Range("M27").Value = 250
$ws.Range("N27").Value = 24000
$ws.Range("O27").Value = 25000
$ws.Range("P27").Value = 24500
$ws.Range("S27").Value = 1225

# Row 28
$ws.Range("D28").Value = 44371
$ws.Range("K28").Value = "Clemenuless"
$ws.Range("L28").Value = "Tercera"
$ws.Range("M28").Value = 120
$ws.Range("N28").Value = 11000
$ws.Range("O28").Value = 12000
$ws.Range("P28").Value = 11500
$ws.Range("R28").Value = "Región Metropolitana"
$ws.Range("S28").Value = 575

# Row 29
$ws.Range("D29").Value = 44342
$ws.Range("M29").Value = 250
$ws.Range("N29").Value = 15000
$ws.Range("O29").Value = 16000
$ws.Range("P29").Value = 15500
$ws.Range("R29").Value = "Región de Coquimbo"
$ws.Range("S29").Value = 775

# Row 30
$ws.Range("D30").Value = 44435
$ws.Range("M30").Value = 300
$ws.Range("N30").Value = 9000
$ws.Range("O30").Value = 10000
$ws.Range("P30").Value = 9500
$ws.Range("S30").Value = 475

# Row 31
$ws.Range("K31").Value = "Murcott"
$ws.Range("L31").Value = "Segunda"
$ws.Range("M31").Value = 250
$ws.Range("N31").Value = 10000
$ws.Range("O31").Value = 11000
$ws.Range("P31").Value = 10500
$ws.Range("S31").Value = 525

# Row 32
$ws.Range("D32").Value = 44231
$ws.Range("N32").Value = 20000
$ws.Range("O32").Value = 21000
$ws.Range("P32").Value = 20500
$ws.Range("S32").Value = 1025

# Row 33
$ws.Range("D33").Value = 44580
$ws.Range("N33").Value = 17000
$ws.Range("O33").Value = 18000
$ws.Range("P33").Value = 17500
$ws.Range("S33").Value = 875
